$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1565
    $ws.Range("F7").Value = 2606
    $ws.Range("F9").Value = 1639
    $ws.Range("F11").Value = 66
    $ws.Range("F12").Value = 534
    $ws.Range("F15").Value = 56
}
